# Reset/update computed TK (temperature coefficient) statistics for each probe row
# so a freshly generated plot selection recalculates with updated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 0.00256
$ws.Range("E4").Value = 122.42
$ws.Range("F4").Value = 0.99874
$ws.Range("G4").Value = 0.00261
$ws.Range("I4").Value = 124.52
$ws.Range("J4").Value = 0.99852

# Row 5
$ws.Range("C5").Value = 0.00334
$ws.Range("E5").Value = 121.76
$ws.Range("F5").Value = 0.99955
$ws.Range("G5").Value = 0.00339
$ws.Range("I5").Value = 124.04
$ws.Range("J5").Value = 0.99922

# Row 6
$ws.Range("C6").Value = 0.0026
$ws.Range("E6").Value = 120.16
$ws.Range("F6").Value = 0.99859
$ws.Range("G6").Value = 0.0027
$ws.Range("I6").Value = 122.69
$ws.Range("J6").Value = 0.9987

# Row 7
$ws.Range("C7").Value = 0.00341
$ws.Range("E7").Value = 120.48
$ws.Range("F7").Value = 0.9996
$ws.Range("G7").Value = 0.00343
$ws.Range("I7").Value = 122.79
$ws.Range("J7").Value = 0.99919

# Row 8
$ws.Range("C8").Value = 0.00235
$ws.Range("E8").Value = 119.9
$ws.Range("F8").Value = 0.99939
$ws.Range("G8").Value = 0.00241
$ws.Range("I8").Value = 122.37
$ws.Range("J8").Value = 0.99843

# Row 9
$ws.Range("C9").Value = 0.00276
$ws.Range("E9").Value = 119.79
$ws.Range("F9").Value = 0.99961
$ws.Range("G9").Value = 0.00278
$ws.Range("I9").Value = 122.18
$ws.Range("J9").Value = 0.9988899999999999

# Row 10
$ws.Range("C10").Value = 0.00263
$ws.Range("E10").Value = 120.06
$ws.Range("F10").Value = 0.99866
$ws.Range("G10").Value = 0.00269
$ws.Range("I10").Value = 122.55
$ws.Range("J10").Value = 0.99866

# Row 11
$ws.Range("C11").Value = 0.00308
$ws.Range("E11").Value = 120.19
$ws.Range("F11").Value = 0.99943
$ws.Range("G11").Value = 0.00314
$ws.Range("I11").Value = 122.39
$ws.Range("J11").Value = 0.99869

# Row 23
$ws.Range("C23").Value = 0.00259
$ws.Range("E23").Value = 120.49
$ws.Range("F23").Value = 0.9968900000000001
$ws.Range("G23").Value = 0.00263
$ws.Range("I23").Value = 122.9
$ws.Range("J23").Value = 0.9983300000000001

# Row 24
$ws.Range("C24").Value = 0.00339
$ws.Range("E24").Value = 120.49
$ws.Range("F24").Value = 0.99955
$ws.Range("G24").Value = 0.00342
$ws.Range("I24").Value = 122.9
$ws.Range("J24").Value = 0.99923

# Row 25
$ws.Range("C25").Value = 0.00258
$ws.Range("E25").Value = 120.49
$ws.Range("F25").Value = 0.99877
$ws.Range("G25").Value = 0.00269
$ws.Range("I25").Value = 122.9
$ws.Range("J25").Value = 0.9987200000000001

# Row 26
$ws.Range("C26").Value = 0.00341
$ws.Range("E26").Value = 120.49
$ws.Range("F26").Value = 0.99918
$ws.Range("I26").Value = 122.9
$ws.Range("J26").Value = 0.9992

# Row 27
$ws.Range("C27").Value = 0.00234
$ws.Range("E27").Value = 120.49
$ws.Range("F27").Value = 0.99898
$ws.Range("G27").Value = 0.00241
$ws.Range("I27").Value = 122.9
$ws.Range("J27").Value = 0.99843

# Row 28
$ws.Range("E28").Value = 120.49
$ws.Range("F28").Value = 0.9996699999999999
$ws.Range("I28").Value = 122.9
$ws.Range("J28").Value = 0.9989

# Row 29
$ws.Range("C29").Value = 0.00261
$ws.Range("E29").Value = 120.49
$ws.Range("F29").Value = 0.9987
$ws.Range("G29").Value = 0.00268
$ws.Range("I29").Value = 122.9
$ws.Range("J29").Value = 0.99866

# Row 30
$ws.Range("C30").Value = 0.00306
$ws.Range("E30").Value = 120.49
$ws.Range("F30").Value = 0.9986
$ws.Range("G30").Value = 0.00312
$ws.Range("I30").Value = 122.9
$ws.Range("J30").Value = 0.99857
